{"js": "// Improved Block error reporting.\n// - Remove the \"Missing start tag of for block :\" paragraph.\n// - Remove the \"End of demonstration.\" paragraph.\n// - In the remaining paragraph, replace the \"A paragraph\" run with the\n//   unbalanced field-code run sequence ({ m:endfor }) plus a bold red\n//   error message run.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs.items;\n\n// Locate the paragraphs by their text content so the script is resilient\n// to ordering assumptions.\nfor (const p of paragraphs) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nlet missingStartPara = null;\nlet endDemoPara = null;\nfor (const p of paragraphs) {\n  const text = p.text;\n  if (text.indexOf(\"Missing start tag of for block\") !== -1) {\n    missingStartPara = p;\n  } else if (text.indexOf(\"End of demonstration.\") !== -1) {\n    endDemoPara = p;\n  }\n}\n\nif (missingStartPara) {\n  missingStartPara.delete();\n}\nif (endDemoPara) {\n  endDemoPara.delete();\n}\nawait context.sync();\n\n// Replace the \"A paragraph\" text with the field-code run sequence and the\n// bold red error message run, while keeping the surrounding bookmark.\nconst found = body.search(\"A paragraph\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  const target = found.items[0];\n\n  const ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    '<w:r><w:instrText xml:space=\"preserve\"> </w:instrText></w:r>' +\n    '<w:r><w:instrText>m:</w:instrText></w:r>' +\n    '<w:r><w:instrText xml:space=\"preserve\">endfor </w:instrText></w:r>' +\n    '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n    '<w:r><w:rPr><w:b w:val=\"true\"/><w:color w:val=\"FF0000\"/></w:rPr>' +\n    '<w:t>Invalid if statement: Unexpected tag m:endfor at this location</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  target.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Improved Block error reporting.\n# - Remove the \"Missing start tag of for block :\" paragraph.\n# - Remove the \"End of demonstration.\" paragraph.\n# - In the remaining paragraph, replace the \"A paragraph\" run with the\n#   unbalanced field-code run sequence ({ m:endfor }) plus a bold red\n#   error message run, while preserving the surrounding bookmark.\n\n$d = $word.ActiveDocument\n\n# --- Replace \"A paragraph\" with the field-code + error-message runs ---\n$rng = $d.Content\n$rng.Find.Text = \"A paragraph\"\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $oldStart = $rng.Start\n    $oldEnd = $rng.End\n\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:fldChar w:fldCharType=\"begin\"/></w:r><w:r><w:instrText xml:space=\"preserve\"> </w:instrText></w:r><w:r><w:instrText>m:</w:instrText></w:r><w:r><w:instrText xml:space=\"preserve\">endfor </w:instrText></w:r><w:r><w:fldChar w:fldCharType=\"end\"/></w:r><w:r><w:rPr><w:b w:val=\"true\"/><w:color w:val=\"FF0000\"/></w:rPr><w:t>Invalid if statement: Unexpected tag m:endfor at this location</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n    # InsertXML inserts the new content right after the matched range;\n    # then remove the original \"A paragraph\" text by its saved position\n    # so the surrounding bookmark (_GoBack) is left untouched.\n    $rng.InsertXML($xml)\n\n    $toDelete = $d.Range($oldStart, $oldEnd)\n    $toDelete.Delete()\n}\n\n# --- Remove the \"Missing start tag of for block :\" and\n#     \"End of demonstration.\" paragraphs (materialize the collection first\n#     since deleting mutates it while iterating) ---\nforeach ($p in @($d.Paragraphs)) {\n    $t = $p.Range.Text\n    if ($t -like \"Missing start tag of for block*\" -or $t -like \"End of demonstration.*\") {\n        $p.Range.Delete()\n    }\n}\n"}
